$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Functions sheet: bump the id column (A2:A9) from 0-based to
#    1-based and move the selection there. Do this BEFORE inserting
#    the new sheet so the new sheet ends up being the active tab.
# ------------------------------------------------------------------
$wsFunctions = $wb.Worksheets.Item("Functions")
$wsFunctions.Activate()
$wsFunctions.Range("A2").Value = 1
$wsFunctions.Range("A3").Value = 2
$wsFunctions.Range("A4").Value = 3
$wsFunctions.Range("A5").Value = 4
$wsFunctions.Range("A6").Value = 5
$wsFunctions.Range("A7").Value = 6
$wsFunctions.Range("A8").Value = 7
$wsFunctions.Range("A9").Value = 8
$wsFunctions.Range("A9").Select()

# ------------------------------------------------------------------
# 2. Insert a new "function_parameters" sheet right after "Users"
#    (i.e. right before "Projects").
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("Users"))
$newSheet.Name = "function_parameters"

# Center (horizontal + vertical) alignment style applied to the whole table.
$tableRange = $newSheet.Range("A1:E7")
$tableRange.HorizontalAlignment = -4108
$tableRange.VerticalAlignment = -4108

# Header row -- write in this specific column order (A,B,C,E,D) so
# that new shared strings are interned in the exact order the
# original workbook used: function_id, type, value, ...
$newSheet.Cells.Item(1,1).Value = "id"
$newSheet.Cells.Item(1,2).Value = "function_id"
$newSheet.Cells.Item(1,3).Value = "kind"
$newSheet.Cells.Item(1,5).Value = "type"
$newSheet.Cells.Item(1,4).Value = "value"

# function_id column (B) -- numeric, order doesn't affect string table.
$newSheet.Cells.Item(2,2).Value = 1
$newSheet.Cells.Item(3,2).Value = 1
$newSheet.Cells.Item(4,2).Value = 2
$newSheet.Cells.Item(5,2).Value = 3
$newSheet.Cells.Item(6,2).Value = 4
$newSheet.Cells.Item(7,2).Value = 5

# kind column (C)
$newSheet.Cells.Item(2,3).Value = "Octopus_Params"
$newSheet.Cells.Item(3,3).Value = "Sys_Params"
$newSheet.Cells.Item(4,3).Value = "text"
$newSheet.Cells.Item(5,3).Value = "Octopus_Params"
$newSheet.Cells.Item(6,3).Value = "Sys_Params"
$newSheet.Cells.Item(7,3).Value = "text"

# value column (D)
$newSheet.Cells.Item(4,4).Value = "ENG"
$newSheet.Cells.Item(7,4).Value = 66

# type column (E)
$newSheet.Cells.Item(2,5).Value = "DataFrame"
$newSheet.Cells.Item(3,5).Value = "DataFrame"
$newSheet.Cells.Item(4,5).Value = "String"
$newSheet.Cells.Item(5,5).Value = "DataFrame"
$newSheet.Cells.Item(6,5).Value = "DataFrame"
$newSheet.Cells.Item(7,5).Value = "string"

# Column widths (characters), closest achievable values to the target
# stored widths 9.625 / 16.75 / 15.875 / 22.25.
$newSheet.Columns.Item(2).ColumnWidth = 8.857142857142858
$newSheet.Columns.Item(3).ColumnWidth = 16
$newSheet.Columns.Item(4).ColumnWidth = 15.142857142857142
$newSheet.Columns.Item(5).ColumnWidth = 21.571428571428573

# Selection on the new sheet.
$newSheet.Range("D2:D3").Select()
